# Add working set of sequences
#
# For every row on Sheet1 whose column F already holds "N/A" (shared string
# index 210), fill columns G through N on that same row with the same
# "N/A" value. This mirrors the "N/A" value already present in column F
# (and C/D) across the rest of the row's data columns (G:N).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToFill = @(7, 12, 15, 16, 21, 23, 26, 31, 32, 38, 40, 45, 46, 48, 51, 54, 58, 65, 66, 67, 69, 72, 76, 79, 82, 88, 92, 94, 96, 98, 99, 103, 107, 108, 110, 111, 115, 119, 120, 123, 128, 131, 135, 136, 138, 139, 141, 145, 152, 155, 159, 160, 161, 163, 167, 171, 174, 175, 176, 178, 179, 182, 187, 193)

foreach ($r in $rowsToFill) {
    $rng = $ws.Range("G" + $r + ":N" + $r)
    $rng.Value = "N/A"
}
